$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4: Desayuno - MILK (2025-08-13 10:27:18)
$ws.Cells.Item(4, 1).Value = "2025-08-13 10:27:18"
$ws.Cells.Item(4, 2).Value = "Desayuno"
$ws.Cells.Item(4, 3).Value = "MILK"
$ws.Cells.Item(4, 4).Value = 1000
$ws.Cells.Item(4, 5).Value = "'620.00"
$ws.Cells.Item(4, 6).Value = "'33.30"
$ws.Cells.Item(4, 7).Value = "'54.20"
$ws.Cells.Item(4, 8).Value = "'33.30"

# New row 5: Almuerzo - BURGER KING, french fries (2025-08-13 10:27:50)
$ws.Cells.Item(5, 1).Value = "2025-08-13 10:27:50"
$ws.Cells.Item(5, 2).Value = "Almuerzo"
$ws.Cells.Item(5, 3).Value = "BURGER KING, french fries"
$ws.Cells.Item(5, 4).Value = 100
$ws.Cells.Item(5, 5).Value = "'280.00"
$ws.Cells.Item(5, 6).Value = "'3.23"
$ws.Cells.Item(5, 7).Value = "'38.70"
$ws.Cells.Item(5, 8).Value = "'0.03"
